# Generate Report for handoff
# Replaces the prior handoff UUID/file info with a new handoff (new source
# file id, new handback xlf hash/timestamps) and inserts a new
# ".localization-config" row describing a failed handoff dependency.

$wb = $excel.ActiveWorkbook

# ---- shared text/values used across sheets -------------------------------
$oldMd        = "1ae56482-212d-48d6-b20a-a338c7c3b051.md"
$newMd        = "dfed2ff0-0c8a-4795-8f89-04c3a8a66ec0.md"
$depMd        = "d482a737-33ce-4a44-8731-1388e99258a8.md"
$cfgName      = ".localization-config"

$notYetHandedOff = "Not yet handed off"
$handoffFailed    = "Handoff failed"
$notLocalized     = "Not localized"

$zhXlf   = "dfed2ff0-0c8a-4795-8f89-04c3a8a66ec0.94c5fc70554ab16089ed18696ace32ea269504c4.zh-cn.xlf"
$deXlf   = "dfed2ff0-0c8a-4795-8f89-04c3a8a66ec0.94c5fc70554ab16089ed18696ace32ea269504c4.de-de.xlf"
$zhTime  = "2016-01-08 12:06:53"
$deTime  = "2016-01-08 12:07:06"
$epoch   = "0001-01-01 00:00:00"
$include = "Include"
$ignored = "Ignored"

$mdTarget    = "https://github.com/OpenLocalizationTest/oltest/blob/c92d59b83f15166c8f7bb2b484e9d8693645cbfd/e2e/$newMd"
$depTarget   = "https://github.com/OpenLocalizationTest/oltest/blob/c92d59b83f15166c8f7bb2b484e9d8693645cbfd/e2e/$depMd"
$cfgTarget   = "https://github.com/OpenLocalizationTest/oltest/blob/c92d59b83f15166c8f7bb2b484e9d8693645cbfd/$cfgName"
$zhXlfTarget = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/da460caece3f43396f1f76d098cfe86d3433fe4e/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/$zhXlf"
$deXlfTarget = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/2cc93161bcee5f8f25f1c73d36fa76102ca2b211/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/$deXlf"

# ===========================================================================
# Sheet "Overview"
# ===========================================================================
$ov = $wb.Worksheets.Item("Overview")
$ov.Hyperlinks.Delete()

$ov.Range("A2").Value = $newMd
$ov.Range("B2").Value = $notYetHandedOff
$ov.Range("C2").Value = $notYetHandedOff

$ov.Range("A3").Value = $depMd
$ov.Range("B3").Value = $handoffFailed
$ov.Range("C3").Value = $handoffFailed

$ov.Range("A4").Value = $cfgName
$ov.Range("B4").Value = $notLocalized
$ov.Range("C4").Value = $notLocalized
$ov.Range("A4").Style = "HyperLink"
$ov.Range("B4").Style = "Normal"
$ov.Range("C4").Style = "Normal"

$ov.Hyperlinks.Add($ov.Range("A2"), $mdTarget, "", "", $newMd)
$ov.Hyperlinks.Add($ov.Range("A3"), $depTarget, "", "", $depMd)
$ov.Hyperlinks.Add($ov.Range("A4"), $cfgTarget, "", "", $cfgName)

# ===========================================================================
# Sheet "zh-cn"
# ===========================================================================
$zh = $wb.Worksheets.Item("zh-cn")
$zh.Hyperlinks.Delete()

$zh.Range("A2").Value = $newMd
$zh.Range("B2").Value = $notYetHandedOff
$zh.Range("C2").Value = $zhXlf
$zh.Range("D2").Value = $zhTime
$zh.Range("D2").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$zh.Range("G2").Value = $epoch
$zh.Range("H2").Value = $include

$zh.Range("A3").Value = $depMd
$zh.Range("B3").Value = $handoffFailed
$zh.Range("D3").Value = $epoch
$zh.Range("D3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$zh.Range("G3").Value = $epoch
$zh.Range("H3").Value = $ignored

$zh.Range("A4").Value = $cfgName
$zh.Range("B4").Value = $notLocalized
$zh.Range("D4").Value = $epoch
$zh.Range("G4").Value = $epoch
$zh.Range("H4").Value = $ignored

$zh.Range("A4").Style = "HyperLink"
$zh.Range("B4").Style = "Normal"
$zh.Range("D4").Style = "Normal"
$zh.Range("D4").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$zh.Range("G4").Style = "Normal"
$zh.Range("H4").Style = "Normal"

$zh.Hyperlinks.Add($zh.Range("A2"), $mdTarget, "", "", $newMd)
$zh.Hyperlinks.Add($zh.Range("C2"), $zhXlfTarget, "", "", $zhXlf)
$zh.Hyperlinks.Add($zh.Range("A3"), $depTarget, "", "", $depMd)
$zh.Hyperlinks.Add($zh.Range("A4"), $cfgTarget, "", "", $cfgName)

# ===========================================================================
# Sheet "de-de"
# ===========================================================================
$de = $wb.Worksheets.Item("de-de")
$de.Hyperlinks.Delete()

$de.Range("A2").Value = $newMd
$de.Range("B2").Value = $notYetHandedOff
$de.Range("C2").Value = $deXlf
$de.Range("D2").Value = $deTime
$de.Range("D2").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$de.Range("G2").Value = $epoch
$de.Range("H2").Value = $include

$de.Range("A3").Value = $depMd
$de.Range("B3").Value = $handoffFailed
$de.Range("D3").Value = $epoch
$de.Range("D3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$de.Range("G3").Value = $epoch
$de.Range("H3").Value = $ignored

$de.Range("A4").Value = $cfgName
$de.Range("B4").Value = $notLocalized
$de.Range("D4").Value = $epoch
$de.Range("G4").Value = $epoch
$de.Range("H4").Value = $ignored

$de.Range("A4").Style = "HyperLink"
$de.Range("B4").Style = "Normal"
$de.Range("D4").Style = "Normal"
$de.Range("D4").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$de.Range("G4").Style = "Normal"
$de.Range("H4").Style = "Normal"

$de.Hyperlinks.Add($de.Range("A2"), $mdTarget, "", "", $newMd)
$de.Hyperlinks.Add($de.Range("C2"), $deXlfTarget, "", "", $deXlf)
$de.Hyperlinks.Add($de.Range("A3"), $depTarget, "", "", $depMd)
$de.Hyperlinks.Add($de.Range("A4"), $cfgTarget, "", "", $cfgName)
